# This script applies updated crafting-leve profit figures (columns H-N)
# to several rows across the ALC, ARM, BSM, CRP, CUL, LTW and WVR sheets
# of the Golem Profits workbook, as produced by the scheduled price-update runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 671.4286
$ws.Range("I2").Value = 501
$ws.Range("J2").Value = 799.25
$ws.Range("K2").Value = 501
$ws.Range("L2").Value = 799.25
$ws.Range("M2").Value = -388
$ws.Range("N2").Value = -1025.25

# Row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 959.4
$ws.Range("I18").Value = 959.4
$ws.Range("K18").Value = 959.4
$ws.Range("M18").Value = -675.4

# Row 105 (Leve Item ID 18668)
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 676.4
$ws.Range("I137").Value = 595.5
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 1786.5
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = 763.5
$ws.Range("N137").Value = -8100

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 6444.2
$ws.Range("J138").Value = 3316
$ws.Range("L138").Value = 9948
$ws.Range("N138").Value = -20228

$ws = $wb.Worksheets.Item("ARM")

# Row 4 (Leve Item ID 5071)
$ws.Range("H4").Value = 127.77778
$ws.Range("I4").Value = 128.375
$ws.Range("K4").Value = 128.375
$ws.Range("M4").Value = -12.375

# Row 95 (Leve Item ID 18204)
$ws.Range("H95").Value = 18800
$ws.Range("J95").Value = 18800
$ws.Range("L95").Value = 18800
$ws.Range("N95").Value = -24292

$ws = $wb.Worksheets.Item("BSM")

# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

# Row 100 (Leve Item ID 18347)
$ws.Range("H100").Value = 85000
$ws.Range("J100").Value = 85000
$ws.Range("L100").Value = 85000
$ws.Range("N100").Value = -87164

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 5559.8
$ws.Range("I134").Value = 5559.8
$ws.Range("K134").Value = 16679.4
$ws.Range("M134").Value = -14144.4

$ws = $wb.Worksheets.Item("CRP")

# Row 15 (Leve Item ID 2406)
$ws.Range("H15").Value = 1878.2174
$ws.Range("J15").Value = 14999.5
$ws.Range("L15").Value = 14999.5
$ws.Range("N15").Value = -15339.5

# Row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 1750

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 1637.8
$ws.Range("I31").Value = 1637.8
$ws.Range("K31").Value = 1637.8
$ws.Range("M31").Value = -1342.8

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 1637.8
$ws.Range("I34").Value = 1637.8
$ws.Range("K34").Value = 1637.8
$ws.Range("M34").Value = -1435.8

# Row 94 (Leve Item ID 32934)
$ws.Range("H94").Value = 1349.5
$ws.Range("J94").Value = 1349.5
$ws.Range("L94").Value = 1349.5
$ws.Range("N94").Value = -2251.5

$ws = $wb.Worksheets.Item("CUL")

# Row 24 (Leve Item ID 4690)
$ws.Range("H24").Value = 673
$ws.Range("J24").Value = 754.5
$ws.Range("L24").Value = 2263.5
$ws.Range("N24").Value = -2723.5

# Row 44 (Leve Item ID 4702)
$ws.Range("H44").Value = 977.7778
$ws.Range("I44").Value = 300
$ws.Range("J44").Value = 1003.8461
$ws.Range("K44").Value = 900
$ws.Range("L44").Value = 3011.5383
$ws.Range("M44").Value = -502
$ws.Range("N44").Value = -3807.5383

# Row 49 (Leve Item ID 4719)
$ws.Range("H49").Value = 4399.4
$ws.Range("J49").Value = 4399.4
$ws.Range("L49").Value = 13198.2
$ws.Range("N49").Value = -13510.2

# Row 87 (Leve Item ID 12864)
$ws.Range("H87").Value = 10690420
$ws.Range("J87").Value = 10690420
$ws.Range("L87").Value = 32071260
$ws.Range("N87").Value = -32073756

# Row 90 (Leve Item ID 12864)
$ws.Range("H90").Value = 10690420
$ws.Range("J90").Value = 10690420
$ws.Range("L90").Value = 96213780
$ws.Range("N90").Value = -96226260

# Row 117 (Leve Item ID 27870)
$ws.Range("H117").Value = 7333.3335
$ws.Range("I117").Value = 1000
$ws.Range("J117").Value = 20000
$ws.Range("K117").Value = 3000
$ws.Range("L117").Value = 60000
$ws.Range("M117").Value = 442
$ws.Range("N117").Value = -66884

$ws = $wb.Worksheets.Item("LTW")

# Row 4 (Leve Item ID 3788)
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 5000
$ws.Range("K4").Value = 5000
$ws.Range("M4").Value = -4887

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 1519.6
$ws.Range("I22").Value = 1098.5
$ws.Range("J22").Value = 1800.3334
$ws.Range("K22").Value = 1098.5
$ws.Range("L22").Value = 1800.3334
$ws.Range("M22").Value = -803.5
$ws.Range("N22").Value = -2390.3334

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 1519.6
$ws.Range("I27").Value = 1098.5
$ws.Range("J27").Value = 1800.3334
$ws.Range("K27").Value = 1098.5
$ws.Range("L27").Value = 1800.3334
$ws.Range("M27").Value = -991.5
$ws.Range("N27").Value = -2014.3334

# Row 28 (Leve Item ID 3788)
$ws.Range("H28").Value = 5000
$ws.Range("I28").Value = 5000
$ws.Range("K28").Value = 5000
$ws.Range("M28").Value = -4768

# Row 37 (Leve Item ID 3788)
$ws.Range("H37").Value = 5000
$ws.Range("I37").Value = 5000
$ws.Range("K37").Value = 5000
$ws.Range("M37").Value = -4893

# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 852208.2
$ws.Range("I40").Value = 18312.25
$ws.Range("K40").Value = 18312.25
$ws.Range("M40").Value = -18176.25

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 497.5
$ws.Range("I46").Value = 497.5
$ws.Range("K46").Value = 497.5
$ws.Range("M46").Value = -309.5

# Row 64 (Leve Item ID 10810)
$ws.Range("H64").Value = 150000
$ws.Range("J64").Value = 150000
$ws.Range("L64").Value = 150000
$ws.Range("N64").Value = -150450

# Row 67 (Leve Item ID 10810)
$ws.Range("H67").Value = 150000
$ws.Range("J67").Value = 150000
$ws.Range("L67").Value = 150000
$ws.Range("N67").Value = -151560

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 716427.4399999999
$ws.Range("I132").Value = 716427.4399999999
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2149282.32
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2146752.32
$ws.Range("N132").ClearContents()

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 626041.1
$ws.Range("I136").Value = 1000959.8
$ws.Range("K136").Value = 3002879.4
$ws.Range("M136").Value = -3000329.4

$ws = $wb.Worksheets.Item("WVR")

# Row 4 (Leve Item ID 2996)
$ws.Range("H4").Value = 1709.9
$ws.Range("J4").Value = 1262.5
$ws.Range("L4").Value = 1262.5
$ws.Range("N4").Value = -1488.5

# Row 18 (Leve Item ID 3543)
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 98 (Leve Item ID 18374)
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 380
$ws.Range("I113").Value = 346.25
$ws.Range("K113").Value = 1038.75
$ws.Range("M113").Value = 1131.25

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1299.4
$ws.Range("I132").Value = 1109.8889
$ws.Range("J132").Value = 3005
$ws.Range("K132").Value = 3329.6667
$ws.Range("L132").Value = 9015
$ws.Range("M132").Value = -799.6666999999998
$ws.Range("N132").Value = -14075

